$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

# Trade #5 closed at 2026-02-16 21:20:49 - leadlag DOWN +0.000%
#
# Seed the new row from the previous trade row (row 4) first. That carries
# over the correctly-typed text cells (Date/Strategy/Side/Status) and the
# still-open Exit Price / Exit Reason placeholders exactly as the rest of
# the table stores them, instead of having Excel's Range.Value auto-detect
# "2026-02-16" as a date literal. Then overwrite just the fields that are
# actually new for trade #5.
$ws.Range("A4:N4").Copy($ws.Range("A5:N5"))

$ws.Range("A5").Value = 5
$ws.Range("C5").Value = "21:20:49"
$ws.Range("F5").Value = 69433.35000000001
$ws.Range("K5").Value = 0.6113
$ws.Range("L5").Value = "Binance leading with -0.061% move"
